# ---------------------------------------------------------------------------
# simple_2.xlsx: insert two new worksheets ("studyDesignArms" and
# "studyDesignEpochs") between "studyDesign" and "mainTimeline", populate
# them with the new Arm / Epoch reference tables, and touch up a couple of
# window/view settings to match the author's edit.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper style sources already present in the workbook -----------------
# studyDesignElements (header row uses style 24, data rows use style 3)
$styleSrc = $wb.Worksheets.Item("studyDesignElements")
# "study" sheet A10 uses the plain vertical-top-only style (style 14)
$vtopSrc  = $wb.Worksheets.Item("study")

# ---------------------------------------------------------------------------
# 1. Create the two new sheets in the right position
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("studyDesign")
$arms = $wb.Worksheets.Add($null, $afterSheet)
$arms.Name = "studyDesignArms"

$epochs = $wb.Worksheets.Add($null, $arms)
$epochs.Name = "studyDesignEpochs"

# ---------------------------------------------------------------------------
# 2. studyDesignArms content
# ---------------------------------------------------------------------------
$arms.Range("A1").Value = "studyArmName"
$arms.Range("B1").Value = "studyArmDescription"
$arms.Range("C1").Value = "studyArmType"
$arms.Range("D1").Value = "studyArmDataOriginDescription"
$arms.Range("E1").Value = "studyArmDataOriginType"

$arms.Range("A2").Value = "Active"
$arms.Range("B2").Value = "Active Substance"
$arms.Range("C2").Value = "Active Comparator Arm"
$arms.Range("D2").Value = "Data collected from subjects"
$arms.Range("E2").Value = "Data Generated Within Study"

$arms.Range("A3").Value = "Placebo"
$arms.Range("B3").Value = "Placebo"
$arms.Range("C3").Value = "Placebo Comparator Arm"
$arms.Range("D3").Value = "Data collected from subjects"
$arms.Range("E3").Value = "Data Generated Within Study"

# Formatting: header row -> style 24 (bold, filled, left/top)
$styleSrc.Range("B1").Copy()
$arms.Range("A1:E1").PasteSpecial(-4122)

# Formatting: data rows -> style 3 (left/top)
$styleSrc.Range("B2").Copy()
$arms.Range("A2:E3").PasteSpecial(-4122)

# Trailing empty formatted cell E4 -> style 14 (vertical-top only)
$vtopSrc.Range("A10").Copy()
$arms.Range("E4").PasteSpecial(-4122)
$arms.Range("E4").ClearContents()

$arms.Columns("A").ColumnWidth = 16.998697916666668
$arms.Columns("B").ColumnWidth = 23.998697916666668
$arms.Columns("C").ColumnWidth = 22.666666666666668
$arms.Columns("D").ColumnWidth = 31.330729166666668
$arms.Columns("E").ColumnWidth = 24.830729166666668

$arms.Activate()
$arms.Range("F29").Select()
$excel.ActiveWindow.Zoom = 150

# ---------------------------------------------------------------------------
# 3. studyDesignEpochs content
# ---------------------------------------------------------------------------
$epochs.Range("A1").Value = "studyEpochName"
$epochs.Range("B1").Value = "studyEpochDescription"
$epochs.Range("C1").Value = "studyEpochType"

$epochs.Range("A2").Value = "Treatment"
$epochs.Range("B2").Value = "Treatment Epoch"
$epochs.Range("C2").Value = "TREATMENT"

$epochs.Range("A3").Value = "Follow-Up"
$epochs.Range("B3").Value = "Follow-up Epoch"
$epochs.Range("C3").Value = "FOLLOW-UP"

$styleSrc.Range("B1").Copy()
$epochs.Range("A1:C1").PasteSpecial(-4122)

$styleSrc.Range("B2").Copy()
$epochs.Range("A2:C3").PasteSpecial(-4122)

$epochs.Columns("A").ColumnWidth = 17.498697916666668
$epochs.Columns("B").ColumnWidth = 26.830729166666668
$epochs.Columns("C").ColumnWidth = 17.498697916666668

$epochs.Activate()
$epochs.Rows("2:2").Select()
$excel.ActiveWindow.Zoom = 160

# ---------------------------------------------------------------------------
# 4. Window / view touch-ups
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.Left = 43440
$win.Top = 4580
